$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: Bitcoin
$ws.Range('D2').NumberFormat = "@"
$ws.Range('D2').Value = '63.831.70'
$ws.Range('E2').Value = '  -1.40%  '

# Row 3: Ethereum
$ws.Range('D3').NumberFormat = "@"
$ws.Range('D3').Value = '3.129.44'

# Row 4: TetherUSD
$ws.Range('E4').Value = '  +0.05%  '

# Row 5: BNB
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '597.25'
$ws.Range('E5').Value = '  -2.61%  '

# Row 6: Solana
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '139.08'
$ws.Range('E6').Value = '  -4.79%  '

# Row 7: USDC
$ws.Range('E7').Value = '  +0.00%  '

# Row 8: LidoStakedEther
$ws.Range('D8').NumberFormat = "@"
$ws.Range('D8').Value = '3.125.95'
$ws.Range('E8').Value = '  -1.12%  '

# Row 9: XRP
$ws.Range('D9').NumberFormat = "@"
$ws.Range('D9').Value = '0.524'
$ws.Range('E9').Value = '  -0.66%  '

# Row 10: Dogecoin
$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '0.147'
$ws.Range('E10').Value = '  -3.37%  '

# Row 11: Toncoin
$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Value = '5.32'
$ws.Range('E11').Value = '  -2.73%  '

# Row 12: Cardano
$ws.Range('D12').NumberFormat = "@"
$ws.Range('D12').Value = '0.462'
$ws.Range('E12').Value = '  -2.54%  '

# Row 13: ShibaInu
$ws.Range('D13').NumberFormat = "@"
$ws.Range('D13').Value = '0.0000251'
$ws.Range('E13').Value = '  -3.13%  '

# Row 14: Avalanche
$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '34.32'
$ws.Range('E14').Value = '  -3.84%  '

# Row 15: WrappedliquidstakedEther2.0
$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '3.644.25'
$ws.Range('E15').Value = '  -1.17%  '

# Row 16: TRON
$ws.Range('E16').Value = '  +2.51%  '

# Row 17: WrappedBTC
$ws.Range('D17').NumberFormat = "@"
$ws.Range('D17').Value = '63.828.29'
$ws.Range('E17').Value = '  -1.35%  '

# Row 18: WrappedEther
$ws.Range('D18').NumberFormat = "@"
$ws.Range('D18').Value = '3.128.81'
$ws.Range('E18').Value = '  -1.17%  '

# Row 19: Polkadot
$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '6.72'
$ws.Range('E19').Value = '  -2.15%  '

# Row 20: BitcoinCash
$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '480.50'
$ws.Range('E20').Value = '  +0.20%  '

# Row 21: Chainlink
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '14.45'
$ws.Range('E21').Value = '  -1.09%  '

# Row 22: Polygon
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '0.702'
$ws.Range('E22').Value = '  -2.68%  '

# Row 23: Uniswap
$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '7.66'
$ws.Range('E23').Value = '  -3.28%  '

# Row 24: Litecoin
$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '87.35'
$ws.Range('E24').Value = '  +3.85%  '

# Row 25: InternetComputer(DFINITY)
$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '13.07'
$ws.Range('E25').Value = '  -4.94%  '

# Row 26: Dai
$ws.Range('E26').Value = '  -0.01%  '

# Row 27: PancakeSwap
$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '2.73'
$ws.Range('E27').Value = '  -2.90%  '

# Row 28: RenderToken
$ws.Range('D28').NumberFormat = "@"
$ws.Range('D28').Value = '8.07'
$ws.Range('E28').Value = '  -7.81%  '

# Row 29: NEARProtocol
$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value = '6.93'
$ws.Range('E29').Value = '  -3.05%  '

# Row 30: ImmutableX
$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value = '2.04'
$ws.Range('E30').Value = '  -3.56%  '

# Row 31: EthereumClassic
$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '27.11'
$ws.Range('E31').Value = '  +2.10%  '

# Row 32: Hedera -> FirstDigitalUSD
$ws.Range('B32').Value = 'FirstDigitalUSD'
$ws.Range('C32').Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '1.00'
$ws.Range('E32').Value = '  +0.03%  '

# Row 33: FirstDigitalUSD -> Hedera
$ws.Range('B33').Value = 'Hedera'
$ws.Range('C33').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '0.111'
$ws.Range('E33').Value = '  -7.80%  '

# Row 34: Stacks
$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '2.58'
$ws.Range('E34').Value = '  -3.79%  '

# Row 35: Mantle
$ws.Range('E35').Value = '  -2.57%  '

# Row 36: Filecoin
$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '5.98'
$ws.Range('E36').Value = '  -0.45%  '

# Row 37: OKB
$ws.Range('D37').NumberFormat = "@"
$ws.Range('D37').Value = '52.57'
$ws.Range('E37').Value = '  -1.38%  '

# Row 38: PEPE
$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '0.0₃0725'
$ws.Range('E38').Value = '  -8.20%  '

# Row 39: VeChain
$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '0.0394'
$ws.Range('E39').Value = '  -1.29%  '

# Row 40: dogwifhat -> Bittensor
$ws.Range('B40').Value = 'Bittensor'
$ws.Range('C40').Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '427.66'
$ws.Range('E40').Value = '  -7.23%  '

# Row 41: Bittensor -> dogwifhat
$ws.Range('B41').Value = 'dogwifhat'
$ws.Range('C41').Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '2.85'
$ws.Range('E41').Value = '  -10.97%  '

# Row 42: Kaspa
$ws.Range('E42').Value = '  -0.79%  '

# Row 43: Cosmos
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '8.25'
$ws.Range('E43').Value = '  -0.92%  '

# Row 44: Maker
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '2.888.94'
$ws.Range('E44').Value = '  +1.03%  '

# Row 45: TheGraph
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '0.256'
$ws.Range('E45').Value = '  -4.18%  '

# Row 46: Fetch.AI -> USDe
$ws.Range('B46').Value = 'USDe'
$ws.Range('C46').Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '0.999'
$ws.Range('E46').Value = '  -0.07%  '

# Row 47: USDe -> Fetch.AI
$ws.Range('B47').Value = 'Fetch.AI'
$ws.Range('C47').Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '2.14'
$ws.Range('E47').Value = '  -8.03%  '

# Row 48: ThetaToken
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '2.35'
$ws.Range('E48').Value = '  -3.76%  '

# Row 49: Stellar
$ws.Range('E49').Value = '  -0.48%  '

# Row 50: InjectiveProtocol
$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '25.44'
$ws.Range('E50').Value = '  -4.42%  '

# Row 51: Monero
$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '120.51'
$ws.Range('E51').Value = '  +0.58%  '
